$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 60 (shifts rows 60-95 down to 61-96),
# representing the new 2020-07-24 data line for Île-de-France / 20-49 salariés.
$ws.Rows.Item(60).Insert()

# Helper to write a value into a cell while forcing "text" storage so that
# numeric-looking strings (counts, amounts) keep their exact textual
# representation (e.g. "30000.00" rather than being coerced to 30000).
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Populate the newly inserted row 60.
Set-TextValue 60 1 "Fonds de solidarité"
Set-TextValue 60 2 "VOLET2"
Set-TextValue 60 3 "3"
Set-TextValue 60 4 "30000.00"
Set-TextValue 60 5 "11"
Set-TextValue 60 6 "Île-de-France"
Set-TextValue 60 7 "12"
Set-TextValue 60 8 "20 à 49 salariés"

# Update the aggregated "nombre_aides" (col C) and "montant_total" (col D)
# values for the rows impacted by the new 2020-07-24 data, using their
# post-insert row numbers.
$updates = @(
    @{Row=39; Col=3; Val="127"},
    @{Row=39; Col=4; Val="296500.00"},
    @{Row=44; Col=3; Val="12"},
    @{Row=44; Col=4; Val="41500.00"},
    @{Row=45; Col=3; Val="46"},
    @{Row=45; Col=4; Val="206122.07"},
    @{Row=46; Col=3; Val="19"},
    @{Row=46; Col=4; Val="104780.00"},
    @{Row=48; Col=3; Val="4"},
    @{Row=48; Col=4; Val="11850.00"},
    @{Row=49; Col=3; Val="80"},
    @{Row=49; Col=4; Val="225937.17"},
    @{Row=50; Col=3; Val="491"},
    @{Row=50; Col=4; Val="1496439.78"},
    @{Row=51; Col=3; Val="211"},
    @{Row=51; Col=4; Val="784686.15"},
    @{Row=52; Col=3; Val="67"},
    @{Row=52; Col=4; Val="354851.23"},
    @{Row=53; Col=3; Val="21"},
    @{Row=53; Col=4; Val="116500.00"},
    @{Row=54; Col=3; Val="14"},
    @{Row=54; Col=4; Val="40720.65"},
    @{Row=55; Col=3; Val="542"},
    @{Row=55; Col=4; Val="1297941.26"},
    @{Row=56; Col=3; Val="2727"},
    @{Row=56; Col=4; Val="7189899.28"},
    @{Row=57; Col=3; Val="1376"},
    @{Row=57; Col=4; Val="4498315.29"},
    @{Row=58; Col=3; Val="471"},
    @{Row=58; Col=4; Val="1892067.00"},
    @{Row=59; Col=3; Val="71"},
    @{Row=59; Col=4; Val="359511.00"},
    @{Row=61; Col=3; Val="205"},
    @{Row=61; Col=4; Val="461423.00"},
    @{Row=74; Col=3; Val="790"},
    @{Row=74; Col=4; Val="2254968.70"},
    @{Row=75; Col=3; Val="296"},
    @{Row=75; Col=4; Val="1073266.79"},
    @{Row=76; Col=3; Val="95"},
    @{Row=76; Col=4; Val="387484.52"},
    @{Row=79; Col=3; Val="192"},
    @{Row=79; Col=4; Val="473026.72"},
    @{Row=80; Col=3; Val="850"},
    @{Row=80; Col=4; Val="2186340.06"},
    @{Row=81; Col=3; Val="301"},
    @{Row=81; Col=4; Val="998675.97"},
    @{Row=82; Col=3; Val="96"},
    @{Row=82; Col=4; Val="384626.30"},
    @{Row=84; Col=3; Val="30"},
    @{Row=84; Col=4; Val="60000.00"}
)

foreach ($u in $updates) {
    Set-TextValue $u.Row $u.Col $u.Val
}
